$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "Belum Ada Laporan Kasus Keracunan MBG, Pemkot Bekasi Tetap Evaluasi SPPGNEWS29/09/2025", "2025-09-29", "Tidak Diketahui", "https://megapolitan.kompas.com/read/2025/09/29/09275851/belum-ada-laporan-kasus-keracunan-mbg-pemkot-bekasi-tetap-evaluasi-sppg", "mbg"),
    @(3, "Usai KLB MBG, Pemerintah Siapkan Ahli Gizi hingga Perbaiki Tata KelolaNEWS29/09/2025", "2025-09-29", "Tidak Diketahui", "https://nasional.kompas.com/read/2025/09/29/08323071/usai-klb-mbg-pemerintah-siapkan-ahli-gizi-hingga-perbaiki-tata-kelola", "mbg"),
    @(4, "Disdik Kalteng Ajak Warga Awasi Program MBG lewat Aplikasi, Bagaimana Caranya?REGIONAL29/09/2025", "2025-09-29", "Tidak Diketahui", "https://regional.kompas.com/read/2025/09/29/083105778/disdik-kalteng-ajak-warga-awasi-program-mbg-lewat-aplikasi-bagaimana", "mbg"),
    @(5, "Kartu Pers Istana Jurnalis CNN Dicabut Usai Tanya soal MBG, Mensesneg Tegaskan Prabowo Tak TerlibatPROV29/09/2025", "2025-09-29", "Tidak Diketahui", "https://www.kompas.com/kalimantan-timur/read/2025/09/29/080000388/kartu-pers-istana-jurnalis-cnn-dicabut-usai-tanya-soal-mbg", "mbg"),
    @(6, "Soal Program MBG, Zulhas: SPPG yang Bermasalah DitutupPROV29/09/2025", "2025-09-29", "Tidak Diketahui", "https://www.kompas.com/jawa-barat/read/2025/09/29/074900188/soal-program-mbg-zulhas--sppg-yang-bermasalah-ditutup", "mbg"),
    @(7, "Duduk Perkara Kartu Pers Istana Jurnalis CNN Dicabut karena Tanya Keracunan MBG ke PrabowoPROV29/09/2025", "2025-09-29", "Tidak Diketahui", "https://www.kompas.com/jawa-tengah/read/2025/09/29/072157488/duduk-perkara-kartu-pers-istana-jurnalis-cnn-dicabut-karena-tanya", "mbg"),
    @(8, "Wali Kota Pastikan Belum Ada Temuan Siswa Keracunan MBG di Banda AcehREGIONAL29/09/2025", "2025-09-29", "Tidak Diketahui", "https://regional.kompas.com/read/2025/09/29/064636878/wali-kota-pastikan-belum-ada-temuan-siswa-keracunan-mbg-di-banda-aceh", "mbg"),
    @(9, "BGN Buka Layanan Pengaduan MBG, Catat Nomor dan Cara Lapornya !MONEY29/09/2025", "2025-09-29", "Tidak Diketahui", "https://money.kompas.com/read/2025/09/29/060125326/bgn-buka-layanan-pengaduan-mbg-catat-nomor-dan-cara-lapornya", "mbg"),
    @(10, "Keracunan MBG, Pers, dan Hari Hak untuk TahuNEWS29/09/2025", "2025-09-29", "Tidak Diketahui", "https://nasional.kompas.com/read/2025/09/29/06000011/keracunan-mbg-pers-dan-hari-hak-untuk-tahu", "mbg"),
    @(11, "Kronologi Pencabutan ID Pers Istana Milik Wartawan CNN Pasca Lempar Pertanyaan Soal MBG ke Presiden PrabowoPROV29/09/2025", "2025-09-29", "Tidak Diketahui", "https://www.kompas.com/jawa-tengah/read/2025/09/29/053000688/kronologi-pencabutan-id-pers-istana-milik-wartawan-cnn-pasca-lempar", "mbg"),
    @(12, "Merasionalkan MBG, Membenahi BGNNEWS29/09/2025", "2025-09-29", "Tidak Diketahui", "https://nasional.kompas.com/read/2025/09/29/05223631/merasionalkan-mbg-membenahi-bgn", "mbg"),
    @(13, "Zulhas soal Keracunan Massal MBG: Ini Bukan Sekadar AngkaMONEY29/09/2025", "2025-09-29", "Tidak Diketahui", "https://money.kompas.com/read/2025/09/29/004600026/zulhas-soal-keracunan-massal-mbg--ini-bukan-sekadar-angka", "mbg")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

# Remove rows 14-18 (they no longer exist in the updated sheet)
$ws.Range("A14:E18").Delete()
